# "fall 22 week 15 day-after inputs"
# Adds a new "Week 49" column (AX) of inning-count data to the InningCounts
# sheet, populating the new week's figures for the players that have data
# so far.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the added week.
$ws.Range("AX1").Value = "Week 49"

# New per-player values for the added week (only some players have
# reported so far).
$ws.Range("AX2").Value = 3.5
$ws.Range("AX4").Value = 5
$ws.Range("AX5").Value = 5.5
$ws.Range("AX7").Value = 4.5
$ws.Range("AX9").Value = 1.5

# Leave the selection on the newly-added last cell, matching where the
# author's cursor ended up after entering the data.
$ws.Range("AX10").Select()
